$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new requirement row above row 13 ("Incluir Vendas") for the new
#    "Monitoramento da produção" requirement, inside the Tabela1 table.
# ---------------------------------------------------------------------------
$ws.Rows(13).Insert()

# Copy formatting from existing cells so the new row matches the table style
# used by the author (B column bold-border style, "normal" border style for
# the TIPO column, and the wrapped/top-aligned style for the description).
$ws.Range("B5").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = "Monitoramento da produção"
$ws.Range("C13").Value = "Funcional"
$ws.Range("D13").Value = "O software tera uma sessão para que seja feita o monitoramento da produção."

# ---------------------------------------------------------------------------
# 2. Update the three hierarchy-attribution descriptions (column G, rows 5-7)
#    to mention the new monitoring access.
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Ter acesso total a área de fornecedores, mas nas áreas do sistema de acesso estoque e acesso vendas apenas visualizar. Acesso ao monitoramento. Não pode ter acesso a área de cadastro de usuários e nem de dar hierarquias."
$ws.Range("G6").Value = "Ter acesso total a área de estoque, mas nas áreas do sistema de acesso fornecedores e acesso vendas apenas visualizar. Acesso ao monitoramento. Não pode ter acesso a área de cadastro de usuários e nem de dar hierarquias."
$ws.Range("G7").Value = "Ter acesso total a área de vendas, mas nas áreas do `nsistema de acesso fornecedores e acesso estoque `napenas visualizar. Acesso ao monitoramento. Não pode ter acesso a área de cadastro de usuários e nem de dar hierarquias."

# ---------------------------------------------------------------------------
# 3. Resize the table / data validation so the new row is included.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:D21"))

$dv = $ws.Range("C4").Validation
$dv.Formula1 = "=`$C`$4:`$C`$21"

# ---------------------------------------------------------------------------
# 4. Row heights — match the re-wrapped heights Excel produced for this edit.
# ---------------------------------------------------------------------------
$ws.Rows(6).RowHeight = 75
$ws.Rows(7).RowHeight = 76.5
$ws.Rows(13).RowHeight = 29.25
$ws.Rows(14).RowHeight = 31.5
$ws.Rows(15).RowHeight = 33
$ws.Rows(16).RowHeight = 30
$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 37.5
$ws.Rows(19).RowHeight = 42.75
$ws.Rows(20).RowHeight = 120
$ws.Rows(21).RowHeight = 60

# ---------------------------------------------------------------------------
# 5. Selection as left by the author after the edit.
# ---------------------------------------------------------------------------
$ws.Range("G11").Select()
